$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the existing row 11, shifting rows 11:13 down to 12:14
$ws.Rows("11:11").Insert()

# Copy the static (unchanging) column values from the row that is now row 12
# into the newly inserted row 11, then set the row-specific values.
$ws.Range("A11").Value = 4
$ws.Range("B11").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C11").Value = "Los Lagos"

$ws.Range("D11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D11").Value = 44778

$ws.Range("E11").Value = 10
$ws.Range("F11").Value = 100112035
$ws.Range("G11").Value = "Bruselas (repollito)"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 120
$ws.Range("K11").Value = 24000
$ws.Range("L11").Value = 24000
$ws.Range("M11").Value = 24000
$ws.Range("N11").Value = "`$/malla 15 kilos"
$ws.Range("O11").Value = "Provincia de Quillota"
$ws.Range("P11").Value = 1600
$ws.Range("Q11").Value = 15
$ws.Range("R11").Value = "Hortaliza"
